# Auto-generated edit script: updates recalculated profit metrics (columns H-N)
# across several Leve tables, per the scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 396.07144
$ws.Range("I53").Value = 429.9
$ws.Range("J53").Value = 377.27777
$ws.Range("K53").Value = 429.9
$ws.Range("L53").Value = 377.27777
$ws.Range("M53").Value = 207.1
$ws.Range("N53").Value = -1651.27777
$ws.Range("H64").Value = 3459.6667
$ws.Range("J64").Value = 3190
$ws.Range("L64").Value = 3190
$ws.Range("N64").Value = -3686
$ws.Range("H67").Value = 3459.6667
$ws.Range("J67").Value = 3190
$ws.Range("L67").Value = 3190
$ws.Range("N67").Value = -4906
$ws.Range("H98").Value = 1956.9166
$ws.Range("I98").Value = 1956.9166
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1956.9166
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -458.9166
$ws.Range("H116").Value = 3734.7036
$ws.Range("I116").Value = 2869.4375
$ws.Range("K116").Value = 2869.4375
$ws.Range("M116").Value = 572.5625
$ws.Range("H122").Value = 1956.9166
$ws.Range("I122").Value = 1956.9166
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5870.7498
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3420.7498
$ws.Range("H132").Value = 6454492
$ws.Range("I132").Value = 7410248.5
$ws.Range("J132").Value = 3132.25
$ws.Range("K132").Value = 22230745.5
$ws.Range("L132").Value = 9396.75
$ws.Range("M132").Value = -22228215.5
$ws.Range("N132").Value = -14456.75
$ws.Range("H137").Value = 4606.8184
$ws.Range("I137").Value = 4380.5557
$ws.Range("K137").Value = 13141.6671
$ws.Range("M137").Value = -10591.6671
$ws.Range("H138").Value = 2763.3867
$ws.Range("I138").Value = 1516.6666
$ws.Range("J138").Value = 4979.778
$ws.Range("K138").Value = 4549.9998
$ws.Range("L138").Value = 14939.334
$ws.Range("M138").Value = 590.0002000000004
$ws.Range("N138").Value = -25219.334
$ws.Range("N98").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 41669584
$ws.Range("I2").Value = 62501124
$ws.Range("K2").Value = 62501124
$ws.Range("M2").Value = -62501011
$ws.Range("H32").Value = 5716.2563
$ws.Range("I32").Value = 4791.2456
$ws.Range("J32").Value = 8227
$ws.Range("K32").Value = 4791.2456
$ws.Range("L32").Value = 8227
$ws.Range("M32").Value = -4504.2456
$ws.Range("N32").Value = -8801
$ws.Range("H45").Value = 1501.4419
$ws.Range("I45").Value = 990.34283
$ws.Range("K45").Value = 990.34283
$ws.Range("M45").Value = -613.34283
$ws.Range("H116").Value = 41669584
$ws.Range("I116").Value = 62501124
$ws.Range("K116").Value = 62501124
$ws.Range("M116").Value = -62498830
$ws.Range("H122").Value = 2630.8823
$ws.Range("I122").Value = 1775.9166
$ws.Range("J122").Value = 4682.8
$ws.Range("K122").Value = 5327.7498
$ws.Range("L122").Value = 14048.4
$ws.Range("M122").Value = -2877.7498
$ws.Range("N122").Value = -18948.4
$ws.Range("H132").Value = 18520836
$ws.Range("I132").Value = 25642514
$ws.Range("J132").Value = 4469.2
$ws.Range("K132").Value = 76927542
$ws.Range("L132").Value = 13407.6
$ws.Range("M132").Value = -76925012
$ws.Range("N132").Value = -18467.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 41669584
$ws.Range("I3").Value = 62501124
$ws.Range("K3").Value = 62501124
$ws.Range("M3").Value = -62501010
$ws.Range("H99").Value = 2491.7273
$ws.Range("I99").Value = 1309.8334
$ws.Range("J99").Value = 3910
$ws.Range("K99").Value = 1309.8334
$ws.Range("L99").Value = 3910
$ws.Range("M99").Value = 188.1666
$ws.Range("N99").Value = -6906
$ws.Range("H105").Value = 2118.8635
$ws.Range("I105").Value = 2080
$ws.Range("J105").Value = 2130.2942
$ws.Range("K105").Value = 2080
$ws.Range("L105").Value = 2130.2942
$ws.Range("M105").Value = -333
$ws.Range("N105").Value = -5624.2942

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2330.7844
$ws.Range("I31").Value = 1811.0638
$ws.Range("J31").Value = 8437.5
$ws.Range("K31").Value = 1811.0638
$ws.Range("L31").Value = 8437.5
$ws.Range("M31").Value = -1516.0638
$ws.Range("N31").Value = -9027.5
$ws.Range("H34").Value = 2330.7844
$ws.Range("I34").Value = 1811.0638
$ws.Range("J34").Value = 8437.5
$ws.Range("K34").Value = 1811.0638
$ws.Range("L34").Value = 8437.5
$ws.Range("M34").Value = -1609.0638
$ws.Range("N34").Value = -8841.5
$ws.Range("H62").Value = 3634.4546
$ws.Range("I62").Value = 2394.75
$ws.Range("K62").Value = 2394.75
$ws.Range("M62").Value = -1770.75
$ws.Range("H65").Value = 3634.4546
$ws.Range("I65").Value = 2394.75
$ws.Range("K65").Value = 11973.75
$ws.Range("M65").Value = -8853.75
$ws.Range("H107").Value = 1413.3334
$ws.Range("I107").Value = 331.16666
$ws.Range("J107").Value = 2495.5
$ws.Range("K107").Value = 331.16666
$ws.Range("L107").Value = 2495.5
$ws.Range("M107").Value = 1588.83334
$ws.Range("N107").Value = -6335.5
$ws.Range("H122").Value = 2318.8667
$ws.Range("I122").Value = 2088.5908
$ws.Range("J122").Value = 2952.125
$ws.Range("K122").Value = 6265.7724
$ws.Range("L122").Value = 8856.375
$ws.Range("M122").Value = -3815.7724
$ws.Range("N122").Value = -13756.375
$ws.Range("H132").Value = 2974.6667
$ws.Range("I132").Value = 2265.8
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 6797.400000000001
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -4267.400000000001
$ws.Range("N132").Value = -20060
$ws.Range("H134").Value = 1684.8462
$ws.Range("I134").Value = 915.35
$ws.Range("J134").Value = 4249.8335
$ws.Range("K134").Value = 2746.05
$ws.Range("L134").Value = 12749.5005
$ws.Range("M134").Value = -211.0500000000002
$ws.Range("N134").Value = -17819.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 11591.777
$ws.Range("I6").Value = 54.333332
$ws.Range("J6").Value = 34666.668
$ws.Range("K6").Value = 162.999996
$ws.Range("L6").Value = 104000.004
$ws.Range("M6").Value = -49.99999600000001
$ws.Range("N6").Value = -104226.004
$ws.Range("H122").Value = 1462.1111
$ws.Range("I122").Value = 517.8570999999999
$ws.Range("J122").Value = 2063
$ws.Range("K122").Value = 4660.7139
$ws.Range("L122").Value = 18567
$ws.Range("M122").Value = -2210.7139
$ws.Range("N122").Value = -23467

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4310.5
$ws.Range("I122").Value = 3305
$ws.Range("J122").Value = 5681.636
$ws.Range("K122").Value = 9915
$ws.Range("L122").Value = 17044.908
$ws.Range("M122").Value = -7465
$ws.Range("N122").Value = -21944.908
$ws.Range("H132").Value = 2625.85
$ws.Range("I132").Value = 2171.3635
$ws.Range("J132").Value = 3875.6875
$ws.Range("K132").Value = 6514.0905
$ws.Range("L132").Value = 11627.0625
$ws.Range("M132").Value = -3984.0905
$ws.Range("N132").Value = -16687.0625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 502724.9
$ws.Range("I122").Value = 590529.3
$ws.Range("J122").Value = 5166.6665
$ws.Range("K122").Value = 1771587.9
$ws.Range("L122").Value = 15499.9995
$ws.Range("M122").Value = -1769137.9
$ws.Range("N122").Value = -20399.9995
$ws.Range("H132").Value = 2861.6882
$ws.Range("I132").Value = 1153.7167
$ws.Range("J132").Value = 8889.823
$ws.Range("K132").Value = 3461.1501
$ws.Range("L132").Value = 26669.469
$ws.Range("M132").Value = -931.1500999999998
$ws.Range("N132").Value = -31729.469
$ws.Range("H136").Value = 786.1045
$ws.Range("I136").Value = 528.6896400000001
$ws.Range("K136").Value = 1586.06892
$ws.Range("M136").Value = 963.9310799999998
